$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.029176170270688
$ws.Cells.Item(2, 3).Value2 = 0.09424712833713755
$ws.Cells.Item(2, 4).Value2 = 0.07680337756706024
$ws.Cells.Item(2, 6).Value2 = 2.957201335550664
$ws.Cells.Item(2, 7).Value2 = 2.21897912280069
$ws.Cells.Item(2, 8).Value2 = 1.844605457834135
$ws.Cells.Item(2, 10).Value2 = 0.2479971667844367
$ws.Cells.Item(2, 11).Value2 = 0.5339279375396586
$ws.Cells.Item(2, 12).Value2 = 0.3109918285135507
$ws.Cells.Item(3, 2).Value2 = 0.9983222857232477
$ws.Cells.Item(3, 3).Value2 = 0.09287175046907947
$ws.Cells.Item(3, 4).Value2 = 0.07616361129435489
$ws.Cells.Item(3, 6).Value2 = 2.952596840799941
$ws.Cells.Item(3, 7).Value2 = 2.212568203605841
$ws.Cells.Item(3, 8).Value2 = 1.846678110595093
$ws.Cells.Item(3, 10).Value2 = 0.2485351462278231
$ws.Cells.Item(3, 11).Value2 = 0.5052137598555362
$ws.Cells.Item(3, 12).Value2 = 0.3077413182180209
$ws.Cells.Item(4, 2).Value2 = 0.9799397599524582
$ws.Cells.Item(4, 3).Value2 = 0.09200715584791297
$ws.Cells.Item(4, 4).Value2 = 0.07580252123761966
$ws.Cells.Item(4, 6).Value2 = 2.950898676154765
$ws.Cells.Item(4, 7).Value2 = 2.209540199742975
$ws.Cells.Item(4, 8).Value2 = 1.848558441881693
$ws.Cells.Item(4, 10).Value2 = 0.2489384194604831
$ws.Cells.Item(4, 11).Value2 = 0.487881853257889
$ws.Cells.Item(4, 12).Value2 = 0.3058886117337565
$ws.Cells.Item(5, 2).Value2 = 0.9725904820591609
$ws.Cells.Item(5, 3).Value2 = 0.09164976188053231
$ws.Cells.Item(5, 4).Value2 = 0.07566339541859435
$ws.Cells.Item(5, 6).Value2 = 2.950490647101887
$ws.Cells.Item(5, 7).Value2 = 2.208534584098601
$ws.Cells.Item(5, 8).Value2 = 1.849477589627725
$ws.Cells.Item(5, 10).Value2 = 0.2491211236943229
$ws.Cells.Item(5, 11).Value2 = 0.4808943659784859
$ws.Cells.Item(5, 12).Value2 = 0.3051697018161263
$ws.Cells.Item(6, 2).Value2 = 0.9713787151113138
$ws.Cells.Item(6, 3).Value2 = 0.0915901107445336
$ws.Cells.Item(6, 4).Value2 = 0.07564077941324854
$ws.Cells.Item(6, 6).Value2 = 2.950440048913265
$ws.Cells.Item(6, 7).Value2 = 2.20838138992211
$ws.Cells.Item(6, 8).Value2 = 1.849639450502025
$ws.Cells.Item(6, 10).Value2 = 0.2491525715956406
$ws.Cells.Item(6, 11).Value2 = 0.4797386620072075
$ws.Cells.Item(6, 12).Value2 = 0.3050525094891157
$ws.Cells.Item(7, 2).Value2 = 0.9798400704633252
$ws.Cells.Item(7, 3).Value2 = 0.09200235642421717
$ws.Cells.Item(7, 4).Value2 = 0.07580061240263802
$ws.Cells.Item(7, 6).Value2 = 2.950892023372631
$ws.Cells.Item(7, 7).Value2 = 2.209525713317717
$ws.Cells.Item(7, 8).Value2 = 1.848570218650892
$ws.Cells.Item(7, 10).Value2 = 0.2489408090758616
$ws.Cells.Item(7, 11).Value2 = 0.4877873118583693
$ws.Cells.Item(7, 12).Value2 = 0.3058787700382481
$ws.Cells.Item(8, 2).Value2 = 1.018421394234991
$ws.Cells.Item(8, 3).Value2 = 0.09377706578501943
$ws.Cells.Item(8, 4).Value2 = 0.07657622599830916
$ws.Cells.Item(8, 6).Value2 = 2.955379423650641
$ws.Cells.Item(8, 7).Value2 = 2.216580056642783
$ws.Cells.Item(8, 8).Value2 = 1.845194024316683
$ws.Cells.Item(8, 10).Value2 = 0.2481675344087861
$ws.Cells.Item(8, 11).Value2 = 0.5239654541192635
$ws.Cells.Item(8, 12).Value2 = 0.3098414026375451
$ws.Cells.Item(9, 2).Value2 = 1.098522673503169
$ws.Cells.Item(9, 3).Value2 = 0.09709823304206822
$ws.Cells.Item(9, 4).Value2 = 0.07834723994054116
$ws.Cells.Item(9, 6).Value2 = 2.973136985386716
$ws.Cells.Item(9, 7).Value2 = 2.237627619089153
$ws.Cells.Item(9, 8).Value2 = 1.843392880378531
$ws.Cells.Item(9, 10).Value2 = 0.2472290935122885
$ws.Cells.Item(9, 11).Value2 = 0.5972724894597832
$ws.Cells.Item(9, 12).Value2 = 0.3187443301003725
$ws.Cells.Item(10, 2).Value2 = 1.160069815274085
$ws.Cells.Item(10, 3).Value2 = 0.0994422228416596
$ws.Cells.Item(10, 4).Value2 = 0.0797987640187543
$ws.Cells.Item(10, 6).Value2 = 2.99164883008693
$ws.Cells.Item(10, 7).Value2 = 2.257502938068853
$ws.Cells.Item(10, 8).Value2 = 1.845006398141095
$ws.Cells.Item(10, 10).Value2 = 0.2468909233522183
$ws.Cells.Item(10, 11).Value2 = 0.6525663123615004
$ws.Cells.Item(10, 12).Value2 = 0.3259721391124373
$ws.Cells.Item(11, 2).Value2 = 1.188652780599995
$ws.Cells.Item(11, 3).Value2 = 0.1004879037252877
$ws.Cells.Item(11, 4).Value2 = 0.08049133038856837
$ws.Cells.Item(11, 6).Value2 = 3.0012584158507
$ws.Cells.Item(11, 7).Value2 = 2.26750594877268
$ws.Cells.Item(11, 8).Value2 = 1.846378004157117
$ws.Cells.Item(11, 10).Value2 = 0.246813158880542
$ws.Cells.Item(11, 11).Value2 = 0.6780320103719362
$ws.Cells.Item(11, 12).Value2 = 0.3294087582529102
$ws.Cells.Item(12, 2).Value2 = 1.199560180414522
$ws.Cells.Item(12, 3).Value2 = 0.1008809240454625
$ws.Cells.Item(12, 4).Value2 = 0.08075818741342999
$ws.Cells.Item(12, 6).Value2 = 3.005068228713426
$ws.Cells.Item(12, 7).Value2 = 2.271432290276636
$ws.Cells.Item(12, 8).Value2 = 1.846989041964378
$ws.Cells.Item(12, 10).Value2 = 0.2467946311778988
$ws.Cells.Item(12, 11).Value2 = 0.6877199292309797
$ws.Cells.Item(12, 12).Value2 = 0.3307314175974909
$ws.Cells.Item(13, 2).Value2 = 1.197207363433279
$ws.Cells.Item(13, 3).Value2 = 0.1007964114487265
$ws.Cells.Item(13, 4).Value2 = 0.08070051118188104
$ws.Cells.Item(13, 6).Value2 = 3.004240118219627
$ws.Cells.Item(13, 7).Value2 = 2.270580523993203
$ws.Cells.Item(13, 8).Value2 = 1.846853369045988
$ws.Cells.Item(13, 10).Value2 = 0.2467981360783043
$ws.Cells.Item(13, 11).Value2 = 0.6856314830430676
$ws.Cells.Item(13, 12).Value2 = 0.3304456139520937
$ws.Cells.Item(14, 2).Value2 = 1.189548465021886
$ws.Cells.Item(14, 3).Value2 = 0.1005202969494334
$ws.Cells.Item(14, 4).Value2 = 0.08051319298160564
$ws.Cells.Item(14, 6).Value2 = 3.001568427084862
$ws.Cells.Item(14, 7).Value2 = 2.267826196352047
$ws.Cells.Item(14, 8).Value2 = 1.846426438261801
$ws.Cells.Item(14, 10).Value2 = 0.2468114158494785
$ws.Cells.Item(14, 11).Value2 = 0.6788281488559562
$ws.Cells.Item(14, 12).Value2 = 0.3295171481367163
$ws.Cells.Item(15, 2).Value2 = 1.184868049136355
$ws.Cells.Item(15, 3).Value2 = 0.1003507840715301
$ws.Cells.Item(15, 4).Value2 = 0.08039905262457125
$ws.Cells.Item(15, 6).Value2 = 2.999954190818471
$ws.Cells.Item(15, 7).Value2 = 2.266157121390449
$ws.Cells.Item(15, 8).Value2 = 1.846176863731671
$ws.Cells.Item(15, 10).Value2 = 0.2468209716428404
$ws.Cells.Item(15, 11).Value2 = 0.674666711410282
$ws.Cells.Item(15, 12).Value2 = 0.3289512054629
$ws.Cells.Item(16, 2).Value2 = 1.158213581506232
$ws.Cells.Item(16, 3).Value2 = 0.09937347102206928
$ws.Cells.Item(16, 4).Value2 = 0.07975414875848941
$ws.Cells.Item(16, 6).Value2 = 2.991044733471668
$ws.Cells.Item(16, 7).Value2 = 2.256868579963168
$ws.Cells.Item(16, 8).Value2 = 1.844929586742609
$ws.Cells.Item(16, 10).Value2 = 0.2468975345660098
$ws.Cells.Item(16, 11).Value2 = 0.6509083276303329
$ws.Cells.Item(16, 12).Value2 = 0.3257505315455376
$ws.Cells.Item(17, 2).Value2 = 1.142011403583751
$ws.Cells.Item(17, 3).Value2 = 0.09876864731899104
$ws.Cells.Item(17, 4).Value2 = 0.07936675396665294
$ws.Cells.Item(17, 6).Value2 = 2.985883442219574
$ws.Cells.Item(17, 7).Value2 = 2.251416754530567
$ws.Cells.Item(17, 8).Value2 = 1.844327695893696
$ws.Cells.Item(17, 10).Value2 = 0.2469639735280396
$ws.Cells.Item(17, 11).Value2 = 0.6364130979513902
$ws.Cells.Item(17, 12).Value2 = 0.323825031124727
$ws.Cells.Item(18, 2).Value2 = 1.132747416024387
$ws.Cells.Item(18, 3).Value2 = 0.09841882817538306
$ws.Cells.Item(18, 4).Value2 = 0.07914697374889101
$ws.Cells.Item(18, 6).Value2 = 2.983026668523905
$ws.Cells.Item(18, 7).Value2 = 2.248371513511188
$ws.Cells.Item(18, 8).Value2 = 1.844041521671812
$ws.Cells.Item(18, 10).Value2 = 0.2470093487280209
$ws.Cells.Item(18, 11).Value2 = 0.628105233841012
$ws.Cells.Item(18, 12).Value2 = 0.3227315314168067
$ws.Cells.Item(19, 2).Value2 = 1.129620262307611
$ws.Cells.Item(19, 3).Value2 = 0.09830005190380575
$ws.Cells.Item(19, 4).Value2 = 0.07907308318032591
$ws.Cells.Item(19, 6).Value2 = 2.982078627726068
$ws.Cells.Item(19, 7).Value2 = 2.247355987330934
$ws.Cells.Item(19, 8).Value2 = 1.843954937511882
$ws.Cells.Item(19, 10).Value2 = 0.247025942523095
$ws.Cells.Item(19, 11).Value2 = 0.6252973936863668
$ws.Cells.Item(19, 12).Value2 = 0.3223636979866029
$ws.Cells.Item(20, 2).Value2 = 1.14373045614019
$ws.Cells.Item(20, 3).Value2 = 0.09883323261877308
$ws.Cells.Item(20, 4).Value2 = 0.07940767858118392
$ws.Cells.Item(20, 6).Value2 = 2.986421293116962
$ws.Cells.Item(20, 7).Value2 = 2.251987742748952
$ws.Cells.Item(20, 8).Value2 = 1.844385556954251
$ws.Cells.Item(20, 10).Value2 = 0.2469561599599643
$ws.Cells.Item(20, 11).Value2 = 0.6379530983345489
$ws.Cells.Item(20, 12).Value2 = 0.324028555926887
$ws.Cells.Item(21, 2).Value2 = 1.191795802443949
$ws.Cells.Item(21, 3).Value2 = 0.1006014786436182
$ws.Cells.Item(21, 4).Value2 = 0.0805680884434139
$ws.Cells.Item(21, 6).Value2 = 3.002348530346836
$ws.Cells.Item(21, 7).Value2 = 2.268631451412801
$ws.Cells.Item(21, 8).Value2 = 1.846549351458094
$ws.Cells.Item(21, 10).Value2 = 0.2468072190422319
$ws.Cells.Item(21, 11).Value2 = 0.6808252445743221
$ws.Cells.Item(21, 12).Value2 = 0.3297892840620023
$ws.Cells.Item(22, 2).Value2 = 1.22369666173222
$ws.Cells.Item(22, 3).Value2 = 0.1017399014700686
$ws.Cells.Item(22, 4).Value2 = 0.08135326001617216
$ws.Cells.Item(22, 6).Value2 = 3.013753857361849
$ws.Cells.Item(22, 7).Value2 = 2.280315949699826
$ws.Cells.Item(22, 8).Value2 = 1.848497645353262
$ws.Cells.Item(22, 10).Value2 = 0.2467735190887268
$ws.Cells.Item(22, 11).Value2 = 0.7091045691662714
$ws.Cells.Item(22, 12).Value2 = 0.3336782843926329
$ws.Cells.Item(23, 2).Value2 = 1.206626125644021
$ws.Cells.Item(23, 3).Value2 = 0.1011338776107138
$ws.Cells.Item(23, 4).Value2 = 0.08093176285376558
$ws.Cells.Item(23, 6).Value2 = 3.007575498284268
$ws.Cells.Item(23, 7).Value2 = 2.274005838723497
$ws.Cells.Item(23, 8).Value2 = 1.847408946681156
$ws.Cells.Item(23, 10).Value2 = 0.2467856883832553
$ws.Cells.Item(23, 11).Value2 = 0.6939876758729042
$ws.Cells.Item(23, 12).Value2 = 0.3315913323103814
$ws.Cells.Item(24, 2).Value2 = 1.142953114021111
$ws.Cells.Item(24, 3).Value2 = 0.09880404014103306
$ws.Cells.Item(24, 4).Value2 = 0.07938916740845769
$ws.Cells.Item(24, 6).Value2 = 2.986177786444259
$ws.Cells.Item(24, 7).Value2 = 2.251729321479075
$ws.Cells.Item(24, 8).Value2 = 1.844359211490314
$ws.Cells.Item(24, 10).Value2 = 0.2469596701122043
$ws.Cells.Item(24, 11).Value2 = 0.6372567842858246
$ws.Cells.Item(24, 12).Value2 = 0.3239365003193342
$ws.Cells.Item(25, 2).Value2 = 1.076378724724435
$ws.Cells.Item(25, 3).Value2 = 0.0962167044619946
$ws.Cells.Item(25, 4).Value2 = 0.07784158077770087
$ws.Cells.Item(25, 6).Value2 = 2.967373645019933
$ws.Cells.Item(25, 7).Value2 = 2.231159958297553
$ws.Cells.Item(25, 8).Value2 = 1.843364309829497
$ws.Cells.Item(25, 10).Value2 = 0.2474212003473042
$ws.Cells.Item(25, 11).Value2 = 0.5771886440400351
$ws.Cells.Item(25, 12).Value2 = 0.3162149952619728